$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 5: the "Từ ngày / {{Start}} / Đến ngày / {{End}}" block moves one
# column to the left (D5:G5 take the old content+format of E5:H5), and
# H5 becomes empty. Capture the original per-cell formatting first (in
# scratch cells far away) before anything gets overwritten, then apply
# both the new text and the carried-over format to each destination.
# ------------------------------------------------------------------
$ws.Range("D5").Copy()
$ws.Range("D200").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E5").Copy()
$ws.Range("E200").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F5").Copy()
$ws.Range("F200").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G5").Copy()
$ws.Range("G200").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H5").Copy()
$ws.Range("H200").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("E200").Copy()
$ws.Range("D5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F200").Copy()
$ws.Range("E5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G200").Copy()
$ws.Range("F5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H200").Copy()
$ws.Range("G5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("D5").Value = "Từ ngày"
$ws.Range("E5").Value = "{{Start}}"
$ws.Range("F5").Value = "Đến ngày"
$ws.Range("G5").Value = "{{End}}"
$ws.Range("H5").Clear()

# clean up the scratch cells used to stage the formats
$ws.Range("D200:H200").Clear()

# ------------------------------------------------------------------
# Row 9: the Stt / StringCreatedAt placeholders move under .Details
# ------------------------------------------------------------------
$ws.Range("A9").Value = "{{ReportStoreStateChange.Details.Stt}}"
$ws.Range("B9").Value = "{{ReportStoreStateChange.Details.StringCreatedAt}}"

# ------------------------------------------------------------------
# Column F gets wider (~17.86 -> ~22.14 characters)
# ------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 21.3

# ------------------------------------------------------------------
# Restore the selection that was active when the sheet was last saved
# ------------------------------------------------------------------
$ws.Range("F12").Select()
